# Natmi following Dr Hou advice
# Update the NATMI LR-pair output (Nts-Sort1) statistics now that the
# ligand/receptor-expressing cell counts changed from 1 to 3 for every
# sending/target cluster combination, which propagates through the
# average/total expression and specificity columns (E..T) for rows 2-9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.905497
$ws.Range("H2").Value = 14.716491
$ws.Range("I2").Value = 0.9446018122065107
$ws.Range("J2").Value = 0.9446018122065107
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.383699
$ws.Range("N2").Value = 4.151097
$ws.Range("O2").Value = 0.08080976933214185
$ws.Range("P2").Value = 0.08080976933214185
$ws.Range("Q2").Value = 6.787731293403
$ws.Range("R2").Value = 61.089581640627
$ws.Range("S2").Value = 0.0763330545551313
$ws.Range("T2").Value = 0.0763330545551313
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.905497
$ws.Range("H3").Value = 14.716491
$ws.Range("I3").Value = 0.9446018122065107
$ws.Range("J3").Value = 0.9446018122065107
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.462094
$ws.Range("N3").Value = 7.386282
$ws.Range("O3").Value = 0.1437893994387872
$ws.Range("P3").Value = 0.1437893994387872
$ws.Range("Q3").Value = 12.077794730718
$ws.Range("R3").Value = 108.700152576462
$ws.Range("S3").Value = 0.1358237272859642
$ws.Range("T3").Value = 0.1358237272859643
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.905497
$ws.Range("H4").Value = 14.716491
$ws.Range("I4").Value = 0.9446018122065107
$ws.Range("J4").Value = 0.9446018122065107
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.336036333333333
$ws.Range("N4").Value = 13.008109
$ws.Range("O4").Value = 0.2532299986575496
$ws.Range("P4").Value = 0.2532299986575496
$ws.Range("Q4").Value = 21.27041322505767
$ws.Range("R4").Value = 191.433719025519
$ws.Range("S4").Value = 0.2392015156369737
$ws.Range("T4").Value = 0.2392015156369737
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.905497
$ws.Range("H5").Value = 14.716491
$ws.Range("I5").Value = 0.9446018122065107
$ws.Range("J5").Value = 0.9446018122065107
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.941088000000001
$ws.Range("N5").Value = 26.823264
$ws.Range("O5").Value = 0.5221708325715213
$ws.Range("P5").Value = 0.5221708325715213
$ws.Range("Q5").Value = 43.860480360736
$ws.Range("R5").Value = 394.744323246624
$ws.Range("S5").Value = 0.4932435147284415
$ws.Range("T5").Value = 0.4932435147284415
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.2876933333333333
$ws.Range("H6").Value = 0.86308
$ws.Range("I6").Value = 0.0553981877934893
$ws.Range("J6").Value = 0.0553981877934893
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.383699
$ws.Range("N6").Value = 4.151097
$ws.Range("O6").Value = 0.08080976933214185
$ws.Range("P6").Value = 0.08080976933214185
$ws.Range("Q6").Value = 0.3980809776399999
$ws.Range("R6").Value = 3.58272879876
$ws.Range("S6").Value = 0.004476714777010547
$ws.Range("T6").Value = 0.004476714777010547
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.2876933333333333
$ws.Range("H7").Value = 0.86308
$ws.Range("I7").Value = 0.0553981877934893
$ws.Range("J7").Value = 0.0553981877934893
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.462094
$ws.Range("N7").Value = 7.386282
$ws.Range("O7").Value = 0.1437893994387872
$ws.Range("P7").Value = 0.1437893994387872
$ws.Range("Q7").Value = 0.7083280298399999
$ws.Range("R7").Value = 6.37495226856
$ws.Range("S7").Value = 0.00796567215282298
$ws.Range("T7").Value = 0.007965672152822982
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2876933333333333
$ws.Range("H8").Value = 0.86308
$ws.Range("I8").Value = 0.0553981877934893
$ws.Range("J8").Value = 0.0553981877934893
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.336036333333333
$ws.Range("N8").Value = 13.008109
$ws.Range("O8").Value = 0.2532299986575496
$ws.Range("P8").Value = 0.2532299986575496
$ws.Range("Q8").Value = 1.247448746191111
$ws.Range("R8").Value = 11.22703871572
$ws.Range("S8").Value = 0.01402848302057598
$ws.Range("T8").Value = 0.01402848302057598
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2876933333333333
$ws.Range("H9").Value = 0.86308
$ws.Range("I9").Value = 0.0553981877934893
$ws.Range("J9").Value = 0.0553981877934893
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.941088000000001
$ws.Range("N9").Value = 26.823264
$ws.Range("O9").Value = 0.5221708325715213
$ws.Range("P9").Value = 0.5221708325715213
$ws.Range("Q9").Value = 2.572291410346667
$ws.Range("R9").Value = 23.15062269312
$ws.Range("S9").Value = 0.0289273178430798
$ws.Range("T9").Value = 0.0289273178430798
